$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3291802
$ws.Range("I40").Value = 10417967
$ws.Range("K40").Value = 10417967
$ws.Range("M40").Value = -10417792
# Row 55
$ws.Range("H55").Value = 1488780.1
$ws.Range("I55").Value = 1154.5454
$ws.Range("J55").Value = 3125168.2
$ws.Range("K55").Value = 1154.5454
$ws.Range("L55").Value = 3125168.2
$ws.Range("M55").Value = -940.5454
$ws.Range("N55").Value = -3125596.2
# Row 138
$ws.Range("H138").Value = 2588.34
$ws.Range("I138").Value = 963.63416
$ws.Range("J138").Value = 3717.3728
$ws.Range("K138").Value = 2890.90248
$ws.Range("L138").Value = 11152.1184
$ws.Range("M138").Value = 2249.09752
$ws.Range("N138").Value = -21432.1184

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 157
$ws.Range("I5").Value = 109.8
$ws.Range("K5").Value = 109.8
$ws.Range("M5").Value = 2.200000000000003
# Row 32
$ws.Range("H32").Value = 38491650
$ws.Range("I32").Value = 58845530
$ws.Range("K32").Value = 58845530
$ws.Range("M32").Value = -58845243
# Row 45
$ws.Range("H45").Value = 27817.29
$ws.Range("I45").Value = 31674.121
$ws.Range("J45").Value = 2362.2
$ws.Range("K45").Value = 31674.121
$ws.Range("L45").Value = 2362.2
$ws.Range("M45").Value = -31297.121
$ws.Range("N45").Value = -3116.2
# Row 110
$ws.Range("H110").Value = 1346.1111
$ws.Range("I110").Value = 711.4545000000001
$ws.Range("K110").Value = 711.4545000000001
$ws.Range("M110").Value = 1333.5455
# Row 122
$ws.Range("H122").Value = 6528.0835
$ws.Range("I122").Value = 8917.25
$ws.Range("K122").Value = 26751.75
$ws.Range("M122").Value = -24301.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 157
$ws.Range("I4").Value = 109.8
$ws.Range("K4").Value = 109.8
$ws.Range("M4").Value = 5.200000000000003
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# Row 105
$ws.Range("H105").Value = 23257786
$ws.Range("I105").Value = 2062
$ws.Range("J105").Value = 83335070
$ws.Range("K105").Value = 2062
$ws.Range("L105").Value = 83335070
$ws.Range("M105").Value = -315
$ws.Range("N105").Value = -83338564
# Row 134
$ws.Range("H134").Value = 1663131.6
$ws.Range("I134").Value = 3430.96
$ws.Range("J134").Value = 6544604
$ws.Range("K134").Value = 10292.88
$ws.Range("L134").Value = 19633812
$ws.Range("M134").Value = -7757.880000000001
$ws.Range("N134").Value = -19638882
# Row 139
$ws.Range("H139").Value = 49000
$ws.Range("J139").Value = 49000
$ws.Range("L139").Value = 49000
$ws.Range("N139").Value = -59280

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 21055.916
$ws.Range("I7").Value = 53.5
$ws.Range("J7").Value = 42058.332
$ws.Range("K7").Value = 53.5
$ws.Range("L7").Value = 42058.332
$ws.Range("M7").Value = 59.5
$ws.Range("N7").Value = -42284.332
# Row 16
$ws.Range("H16").Value = 3600
$ws.Range("I16").Value = 2800
$ws.Range("J16").Value = 3760
$ws.Range("K16").Value = 2800
$ws.Range("L16").Value = 3760
$ws.Range("M16").Value = -2513
$ws.Range("N16").Value = -4334
# Row 31
$ws.Range("H31").Value = 3553.2812
$ws.Range("I31").Value = 1991.0769
$ws.Range("J31").Value = 4622.1577
$ws.Range("K31").Value = 1991.0769
$ws.Range("L31").Value = 4622.1577
$ws.Range("M31").Value = -1696.0769
$ws.Range("N31").Value = -5212.1577
# Row 34
$ws.Range("H34").Value = 3553.2812
$ws.Range("I34").Value = 1991.0769
$ws.Range("J34").Value = 4622.1577
$ws.Range("K34").Value = 1991.0769
$ws.Range("L34").Value = 4622.1577
$ws.Range("M34").Value = -1789.0769
$ws.Range("N34").Value = -5026.1577
# Row 68
$ws.Range("H68").Value = 22795
$ws.Range("J68").Value = 22795
$ws.Range("L68").Value = 22795
$ws.Range("N68").Value = -24293
# Row 71
$ws.Range("H71").Value = 22795
$ws.Range("J71").Value = 22795
$ws.Range("L71").Value = 68385
$ws.Range("N71").Value = -75873
# Row 80
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -32246
# Row 83
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -101232
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 99
$ws.Range("H99").Value = 90920670
$ws.Range("I99").Value = 125014296
$ws.Range("J99").Value = 4333.3335
$ws.Range("K99").Value = 125014296
$ws.Range("L99").Value = 4333.3335
$ws.Range("M99").Value = -125012798
$ws.Range("N99").Value = -7329.3335
# Row 113
$ws.Range("H113").Value = 3600
$ws.Range("I113").Value = 2800
$ws.Range("J113").Value = 3760
$ws.Range("K113").Value = 2800
$ws.Range("L113").Value = 3760
$ws.Range("M113").Value = -630
$ws.Range("N113").Value = -8100
# Row 122
$ws.Range("H122").Value = 83334330
$ws.Range("I122").Value = 83334330
$ws.Range("K122").Value = 250002990
$ws.Range("M122").Value = -250000540
# Row 126
$ws.Range("H126").Value = 90920670
$ws.Range("I126").Value = 125014296
$ws.Range("J126").Value = 4333.3335
$ws.Range("K126").Value = 375042888
$ws.Range("L126").Value = 13000.0005
$ws.Range("M126").Value = -375040418
$ws.Range("N126").Value = -17940.0005
# Row 134
$ws.Range("H134").Value = 2092.2424
$ws.Range("I134").Value = 2349.8096
$ws.Range("J134").Value = 1641.5
$ws.Range("K134").Value = 7049.4288
$ws.Range("L134").Value = 4924.5
$ws.Range("M134").Value = -4514.4288
$ws.Range("N134").Value = -9994.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 3459.158
$ws.Range("I3").Value = 2734
$ws.Range("J3").Value = 7326.6665
$ws.Range("K3").Value = 8202
$ws.Range("L3").Value = 21979.9995
$ws.Range("M3").Value = -8090
$ws.Range("N3").Value = -22203.9995
# Row 113
$ws.Range("H113").Value = 26191106
$ws.Range("J113").Value = 20000720
$ws.Range("L113").Value = 60002160
$ws.Range("N113").Value = -60006500
# Row 123
$ws.Range("H123").Value = 203.33333
$ws.Range("I123").Value = 203.33333
$ws.Range("K123").Value = 609.99999
$ws.Range("M123").Value = 1840.00001
# Row 131
$ws.Range("H131").Value = 755.25
$ws.Range("I131").Value = 525
$ws.Range("J131").Value = 769.94684
$ws.Range("K131").Value = 1575
$ws.Range("L131").Value = 2309.84052
$ws.Range("M131").Value = 3465
$ws.Range("N131").Value = -12389.84052

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 39800
$ws.Range("J15").Value = 39800
$ws.Range("L15").Value = 39800
$ws.Range("N15").Value = -40376
# Row 57
$ws.Range("H57").Value = 15030.5
$ws.Range("J57").Value = 15030.5
$ws.Range("L57").Value = 15030.5
$ws.Range("N57").Value = -16670.5
# Row 81
$ws.Range("H81").Value = 39800
$ws.Range("J81").Value = 39800
$ws.Range("L81").Value = 39800
$ws.Range("N81").Value = -41796
# Row 84
$ws.Range("H84").Value = 39800
$ws.Range("J84").Value = 39800
$ws.Range("L84").Value = 119400
$ws.Range("N84").Value = -129384
# Row 102
$ws.Range("H102").Value = 1365.7142
$ws.Range("I102").Value = 1198
$ws.Range("J102").Value = 1785
$ws.Range("K102").Value = 1198
$ws.Range("L102").Value = 1785
$ws.Range("M102").Value = 424
$ws.Range("N102").Value = -5029
# Row 113
$ws.Range("H113").Value = 2660
$ws.Range("I113").Value = 2833.3333
$ws.Range("J113").Value = 2400
$ws.Range("K113").Value = 2833.3333
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = -663.3332999999998
$ws.Range("N113").Value = -6740
# Row 132
$ws.Range("H132").Value = 6900.3213
$ws.Range("I132").Value = 2025.3334
$ws.Range("J132").Value = 12525.308
$ws.Range("K132").Value = 6076.0002
$ws.Range("L132").Value = 37575.924
$ws.Range("M132").Value = -3546.0002
$ws.Range("N132").Value = -42635.924

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 10417341
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 1000
$ws.Range("N46").Value = -1376
# Row 61
$ws.Range("H61").Value = 1851.95
$ws.Range("I61").Value = 1556.6666
$ws.Range("J61").Value = 2737.8
$ws.Range("K61").Value = 1556.6666
$ws.Range("L61").Value = 2737.8
$ws.Range("M61").Value = -1354.6666
$ws.Range("N61").Value = -3141.8
# Row 113
$ws.Range("H113").Value = 1851.95
$ws.Range("I113").Value = 1556.6666
$ws.Range("J113").Value = 2737.8
$ws.Range("K113").Value = 1556.6666
$ws.Range("L113").Value = 2737.8
$ws.Range("M113").Value = 613.3334
$ws.Range("N113").Value = -7077.8
# Row 122
$ws.Range("H122").Value = 80000
$ws.Range("I122").Value = 80000
$ws.Range("K122").Value = 240000
$ws.Range("M122").Value = -237550

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3675.3333
$ws.Range("I122").Value = 2701
$ws.Range("K122").Value = 8103
$ws.Range("M122").Value = -5653
# Row 126
$ws.Range("H126").Value = 1607.6154
$ws.Range("I126").Value = 1081.7273
$ws.Range("K126").Value = 3245.1819
$ws.Range("M126").Value = -775.1819
# Row 136
$ws.Range("H136").Value = 5449.6304
$ws.Range("I136").Value = 9746.143
$ws.Range("J136").Value = 1840.56
$ws.Range("K136").Value = 29238.429
$ws.Range("L136").Value = 5521.68
$ws.Range("M136").Value = -26688.429
$ws.Range("N136").Value = -10621.68
# Row 139
$ws.Range("H139").Value = 49367
$ws.Range("J139").Value = 49326.875
$ws.Range("L139").Value = 49326.875
$ws.Range("N139").Value = -59606.875

Write-Host "All updates applied."
